$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new cell values look like plain numbers (e.g. "604.97").
# The source data stores them as text, so if we just set .Value, Excel
# will auto-convert the string into a real number (and mangle trailing
# zeros, e.g. "8.60" -> 8.6). To keep them as text - like the original
# cells - we temporarily force Text number format on those specific
# cells, assign the values, then restore the default (Normal) style so
# the saved file keeps the same (unstyled) cells as before.
$numericTextCells = "D5","D6","D11","D15","D16","D19","D21","D22","D23","D24","D25","D26","D28","D29","D30","D31","D32","D34","D37","D38","D40","D46","D48","D51"
$numericTextRange = $ws.Range($numericTextCells -join ",")
foreach ($area in $numericTextRange.Areas) {
    $area.NumberFormat = "@"
}

$ws.Range("D2").Value = "70.310.89"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "3.609.97"
$ws.Range("E3").Value = "  +1.63%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "604.97"
$ws.Range("D6").Value = "195.72"
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("D11").Value = "53.62"
$ws.Range("E11").Value = "  -1.07%  "
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").Value = "4.184.14"
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("D15").Value = "12.99"
$ws.Range("E15").Value = "  +1.95%  "
$ws.Range("D16").Value = "595.82"
$ws.Range("E16").Value = "  -1.59%  "
$ws.Range("D17").Value = "70.426.73"
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.614.87"
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "19.02"
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").Value = "0.996"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").Value = "17.84"
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("D23").Value = "5.18"
$ws.Range("E23").Value = "  -2.26%  "
$ws.Range("D24").Value = "101.74"
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("D25").Value = "4.63"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "3.02"
$ws.Range("E26").Value = "  -3.83%  "
$ws.Range("E27").Value = "  -1.69%  "
$ws.Range("D28").Value = "9.59"
$ws.Range("E28").Value = "  -0.86%  "
$ws.Range("D29").Value = "33.73"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").Value = "4.70"
$ws.Range("E30").Value = "  +4.24%  "
$ws.Range("D31").Value = "7.20"
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("D32").Value = "12.31"
$ws.Range("E32").Value = "  -3.16%  "
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("D34").Value = "63.55"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").Value = "0.0₃0894"
$ws.Range("E35").Value = "  +6.36%  "
$ws.Range("D36").Value = "3.898.82"
$ws.Range("E36").Value = "  +3.05%  "
$ws.Range("D37").Value = "542.51"
$ws.Range("E37").Value = "  +11.33%  "
$ws.Range("D38").Value = "3.13"
$ws.Range("E38").Value = "  +1.22%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").Value = "36.91"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("E41").Value = "  -1.58%  "
$ws.Range("E42").Value = "  -4.21%  "
$ws.Range("E43").Value = "  -1.67%  "
$ws.Range("E44").Value = "  -0.65%  "
$ws.Range("E45").Value = "  +4.14%  "
$ws.Range("D46").Value = "2.86"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("D48").Value = "8.60"
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("D51").Value = "1.30"
$ws.Range("E51").Value = "  -0.30%  "

foreach ($area in $numericTextRange.Areas) {
    $area.Style = "Normal"
}

